# Updated symbol list on Wed Jan 11 23:26:50 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the crypto
# listing on the active sheet. Values are kept as plain text (matching the
# original string cells) rather than being auto-converted by Excel into
# numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "280.90";     E = "1.01%" }
    @{ Row = 3;  D = "27.56";      E = "1.35%" }
    @{ Row = 4;  D = "4.854";      E = "-1.15%" }
    @{ Row = 5;  D = "0.06439";    E = "0.39%" }
    @{ Row = 6;  D = "7.102";      E = "2.14%" }
    @{ Row = 7;  D = "1.279";      E = "2.73%" }
    @{ Row = 8;  D = "0.9020";     E = "2.33%" }
    @{ Row = 9;  D = "0.1546";     E = "1.38%" }
    @{ Row = 10; D = "0.06323";    E = "25.71%" }
    @{ Row = 11; D = "0.07503";    E = "-0.21%" }
    @{ Row = 12; D = "0.02927";    E = "1.61%" }
    @{ Row = 13; D = "0.08988";    E = "-0.22%" }
    @{ Row = 14; D = "0.001576";   E = "-0.40%" }
    @{ Row = 15; D = "0.0006408";  E = "-0.22%" }
    @{ Row = 16; D = "0.006033";   E = "2.92%" }
    @{ Row = 18; D = "3.304" }
    @{ Row = 19; D = "2.234";      E = "-1.69%" }
    @{ Row = 21; E = "1.16%" }
    @{ Row = 22; D = "3.900";      E = "-0.06%" }
    @{ Row = 23; D = "0.04415";    E = "-0.29%" }
    @{ Row = 24; E = "8.79%" }
    @{ Row = 25; D = "0.001175";   E = "0.09%" }
    @{ Row = 26; D = "0.004307";   E = "11.55%" }
    @{ Row = 28; D = "0.0001179";  E = "-1.79%" }
    @{ Row = 29; D = "0.0001655" }
    @{ Row = 40; D = "0.04092";    E = "-1.18%" }
    @{ Row = 41; D = "0.1403";     E = "19.38%" }
    @{ Row = 42; D = "0.006614";   E = "-3.02%" }
    @{ Row = 43; D = "0.002068";   E = "-13.58%" }
    @{ Row = 44; D = "0.01161";    E = "-1.02%" }
    @{ Row = 45; D = "0.00005562"; E = "6.87%" }
    @{ Row = 46; D = "1.628";      E = "9.53%" }
    @{ Row = 47; E = "-8.86%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        # Leading apostrophe keeps Excel from re-interpreting the text as a
        # number; resetting the style afterwards avoids leaving a stray
        # quote-prefix/number-format on the cell so it stays identical to
        # its original (unstyled) look.
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.Value = "'" + $u.E
        $cell.Style = "Normal"
    }
}

$wb.Save()
